# Category_Level_Intensity_Analysis.xlsx — re-upload edit
#
# The interval-label column (column A, rows 3-10 and 12-19) on every
# Gain_*/Loss_* sheet previously held a distinct "YYYY-YYYY" string per
# row (1996-2008 .. 1996-2015, 2007-2017 .. 2007-2024). Those are
# collapsed back down to the two canonical interval labels used by row 2
# / row 11 ("1996-2007" and "2007-2016"), which removes the now-unused
# unique strings from the shared string table on save.
#
# On the Loss_* sheets, the A1 header cell ("Interval*") also picks up
# the same font styling used on the Gain_* sheets' A1 header (nudging
# Font.Name/Font.Size onto their current values forces the style index
# to be resolved fresh instead of keeping the stale one).

$wb = $excel.ActiveWorkbook

$sheetNames = @(
    "Gain_AYE", "Loss_AYE",
    "Gain_BAG", "Loss_BAG",
    "Gain_MON", "Loss_MON",
    "Gain_RAK", "Loss_RAK",
    "Gain_TNI", "Loss_TNI",
    "Gain_YGN", "Loss_YGN"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    for ($r = 3; $r -le 10; $r++) {
        $ws.Cells.Item($r, 1).Value = "1996-2007"
    }
    for ($r = 12; $r -le 19; $r++) {
        $ws.Cells.Item($r, 1).Value = "2007-2016"
    }

    if ($name.StartsWith("Loss_")) {
        $a1 = $ws.Range("A1")
        $a1.Font.Name = "Calibri"
        $a1.Font.Size = 12
    }
}
